$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.432.32"
$ws.Range("E2").Value = "  -1.50%  "

$ws.Range("D3").Value = "'1.892.81"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("D4").Value = "'0.9973"
$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").Value = "'237.60"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("D7").Value = "'0.4840"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("D8").Value = "'0.2898"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("D9").Value = "'0.06616"
$ws.Range("E9").Value = "  -2.30%  "

$ws.Range("D10").Value = "'1.911.71"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").Value = "'16.96"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").Value = "'0.07344"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "'5.162"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").Value = "'88.00"
$ws.Range("E14").Value = "  -1.67%  "

$ws.Range("D15").Value = "'0.6641"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").Value = "'30.408.50"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").Value = "'13.43"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "'0.000007792"
$ws.Range("E18").Value = "  -2.53%  "

$ws.Range("D19").Value = "'0.9973"
$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").Value = "'5.419"
$ws.Range("E20").Value = "  +4.79%  "

$ws.Range("D21").Value = "'2.121.88"
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("D22").Value = "'0.9968"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("D23").Value = "'196.38"
$ws.Range("E23").Value = "  -5.00%  "

$ws.Range("D24").Value = "'6.211"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").Value = "'9.335"
$ws.Range("E25").Value = "  -3.48%  "

$ws.Range("D26").Value = "'165.23"
$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("D27").Value = "'18.21"
$ws.Range("E27").Value = "  -4.35%  "

$ws.Range("D28").Value = "'1.942"
$ws.Range("E28").Value = "  -2.27%  "

$ws.Range("D29").Value = "'1.443"
$ws.Range("E29").Value = "  +1.08%  "

$ws.Range("D30").Value = "'4.332"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").Value = "'0.09170"
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("D32").Value = "'4.042"
$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("D33").Value = "'0.05101"
$ws.Range("E33").Value = "  -1.74%  "

$ws.Range("D34").Value = "'1.157"
$ws.Range("E34").Value = "  +3.15%  "

$ws.Range("D35").Value = "'0.7317"
$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("D36").Value = "'2.691"
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("E37").Value = "  -3.76%  "

$ws.Range("D38").Value = "'2.641"
$ws.Range("E38").Value = "  -3.29%  "

$ws.Range("D39").Value = "'0.9223"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").Value = "'2.078"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").Value = "'106.36"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").Value = "'0.4331"
$ws.Range("E42").Value = "  -3.93%  "

$ws.Range("D43").Value = "'5.849"
$ws.Range("E43").Value = "  -1.32%  "

$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("D45").Value = "'7.571"
$ws.Range("E45").Value = "  -1.56%  "

$ws.Range("D46").Value = "'0.1323"
$ws.Range("E46").Value = "  -5.15%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'65.44"
$ws.Range("E47").Value = "  -9.44%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.537"
$ws.Range("E48").Value = "  +6.73%  "

$ws.Range("D49").Value = "'8.948"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("D50").Value = "'34.04"
$ws.Range("E50").Value = "  -5.19%  "

$ws.Range("D51").Value = "'0.05756"
$ws.Range("E51").Value = "  -3.27%  "
